# Update "想去人数" (want-to-go count, column F) values across the
# 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4) sheets, matching the
# gh-pages data refresh captured in the commit's generated output.

$wb = $excel.ActiveWorkbook

# --- 展览 (Worksheet 1) -----------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 163
$ws1.Range("F3").Value  = 286
$ws1.Range("F4").Value  = 123
$ws1.Range("F5").Value  = 1266
$ws1.Range("F6").Value  = 17776
$ws1.Range("F7").Value  = 343
$ws1.Range("F8").Value  = 244
$ws1.Range("F10").Value = 6716
$ws1.Range("F13").Value = 8
$ws1.Range("F15").Value = 56
$ws1.Range("F16").Value = 5
$ws1.Range("F18").Value = 1296
$ws1.Range("F19").Value = 182
$ws1.Range("F25").Value = 259
$ws1.Range("F26").Value = 964
$ws1.Range("F27").Value = 105
$ws1.Range("F28").Value = 5135
$ws1.Range("F31").Value = 11902
$ws1.Range("F32").Value = 1267
$ws1.Range("F35").Value = 261
$ws1.Range("F36").Value = 3905
$ws1.Range("F37").Value = 290
$ws1.Range("F38").Value = 89

# --- 演出 (Worksheet 2) -------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value  = 26

# --- 全部类型 (Worksheet 4) ---------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 163
$ws4.Range("F3").Value  = 286
$ws4.Range("F4").Value  = 123
$ws4.Range("F5").Value  = 1266
$ws4.Range("F6").Value  = 17776
$ws4.Range("F7").Value  = 343
$ws4.Range("F8").Value  = 244
$ws4.Range("F10").Value = 6716
$ws4.Range("F13").Value = 8
$ws4.Range("F15").Value = 56
$ws4.Range("F16").Value = 5
$ws4.Range("F18").Value = 1296
$ws4.Range("F19").Value = 182
$ws4.Range("F25").Value = 259
$ws4.Range("F26").Value = 964
$ws4.Range("F27").Value = 105
$ws4.Range("F28").Value = 5135
$ws4.Range("F31").Value = 26
$ws4.Range("F33").Value = 11902
$ws4.Range("F34").Value = 1267
$ws4.Range("F37").Value = 261
$ws4.Range("F38").Value = 3905
$ws4.Range("F39").Value = 290
$ws4.Range("F40").Value = 89
